$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "248.62"
    "D3" = "21.74"
    "D4" = "5.289"
    "D5" = "0.05595"
    "D7" = "6.374"
    "D8" = "0.8152"
    "D9" = "0.9741"
    "D10" = "0.1416"
    "D11" = "0.07572"
    "D12" = "0.03140"
    "D13" = "0.03042"
    "D14" = "0.09306"
    "D15" = "3.567"
    "D16" = "0.001601"
    "D17" = "0.04686"
    "D18" = "0.0005770"
    "D19" = "0.006464"
    "D20" = "0.005037"
    "D21" = "0.001032"
    "D23" = "3.748"
    "D24" = "2.142"
    "D40" = "0.03942"
    "D41" = "0.007025"
    "D42" = "0.1050"
    "D43" = "0.003394"
    "D44" = "0.008514"
    "D45" = "0.00005807"
    "D47" = "0.0005490"
    "D48" = "0.6788"
    "D49" = "0.1546"
    "D50" = "0.00002096"
    "D51" = "0.01008"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.Value = "'" + $updates[$ref]
    $cell.Style = $origStyle
}
